$d = $word.ActiveDocument

# --- Insert a bold paragraph "Text vor der Tabelle" before the table ---
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("Text vor der Tabelle`r")
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.Font.Bold = 1

# --- Insert an italic run "Text nach der Tabelle" into the paragraph that
#     follows the table (the one holding the _GoBack bookmark), right
#     before the bookmark, without creating a new paragraph ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertPos = $lastPara.Range.Start
$insertRng = $d.Range($insertPos, $insertPos)
$insertRng.InsertBefore("Text nach der Tabelle")

$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.Font.Italic = 1

Write-Output "done"
